$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 308, shifting existing data down.
$ws.Rows("308:310").Insert()

# Common constant values shared by every record in this data block.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$categoriaId = 100112002
$categoria = "Pimiento"
$clasificacion = "Hortaliza"

# New row 308: Cuatro cascos rojo
$ws.Cells.Item(308, 1).Value = $mercadoId
$ws.Cells.Item(308, 2).Value = $mercado
$ws.Cells.Item(308, 3).Value = $region
$ws.Cells.Item(308, 4).Value = 44964
$ws.Cells.Item(308, 5).Value = $codreg
$ws.Cells.Item(308, 6).Value = $categoriaId
$ws.Cells.Item(308, 7).Value = $categoria
$ws.Cells.Item(308, 8).Value = "Cuatro cascos rojo"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 50
$ws.Cells.Item(308, 11).Value = 14000
$ws.Cells.Item(308, 12).Value = 14000
$ws.Cells.Item(308, 13).Value = 14000
$ws.Cells.Item(308, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(308, 15).Value = "Región del Maule"
$ws.Cells.Item(308, 16).Value = 778
$ws.Cells.Item(308, 17).Value = 18
$ws.Cells.Item(308, 18).Value = $clasificacion

# New row 309: Cuatro cascos verde
$ws.Cells.Item(309, 1).Value = $mercadoId
$ws.Cells.Item(309, 2).Value = $mercado
$ws.Cells.Item(309, 3).Value = $region
$ws.Cells.Item(309, 4).Value = 44964
$ws.Cells.Item(309, 5).Value = $codreg
$ws.Cells.Item(309, 6).Value = $categoriaId
$ws.Cells.Item(309, 7).Value = $categoria
$ws.Cells.Item(309, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 60
$ws.Cells.Item(309, 11).Value = 7000
$ws.Cells.Item(309, 12).Value = 7500
$ws.Cells.Item(309, 13).Value = 7250
$ws.Cells.Item(309, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(309, 15).Value = "Región del Maule"
$ws.Cells.Item(309, 16).Value = 403
$ws.Cells.Item(309, 17).Value = 18
$ws.Cells.Item(309, 18).Value = $clasificacion

# New row 310: Zafiro rojo
$ws.Cells.Item(310, 1).Value = $mercadoId
$ws.Cells.Item(310, 2).Value = $mercado
$ws.Cells.Item(310, 3).Value = $region
$ws.Cells.Item(310, 4).Value = 44964
$ws.Cells.Item(310, 5).Value = $codreg
$ws.Cells.Item(310, 6).Value = $categoriaId
$ws.Cells.Item(310, 7).Value = $categoria
$ws.Cells.Item(310, 8).Value = "Zafiro rojo"
$ws.Cells.Item(310, 9).Value = "Primera"
$ws.Cells.Item(310, 10).Value = 60
$ws.Cells.Item(310, 11).Value = 18000
$ws.Cells.Item(310, 12).Value = 18000
$ws.Cells.Item(310, 13).Value = 18000
$ws.Cells.Item(310, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(310, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(310, 16).Value = 1200
$ws.Cells.Item(310, 17).Value = 15
$ws.Cells.Item(310, 18).Value = $clasificacion

Write-Output "done"
